# Script execution via Springboot configured and working
#
# TC_Details originally had 5 columns:
#   A=TestCaseID, B=TestCaseName, C=Execute_Flag, D=AutomationScriptName, E=Remarks
# with two data rows (GoogleLaunchURLTest / ManulifeLaunchURLTest).
#
# This edit:
#   1. Removes the TestCaseID / TestCaseName columns (A & B) -- everything
#      shifts left so Execute_Flag, AutomationScriptName and Remarks become
#      columns A, B and C.
#   2. Flips the second test row's flag/remark from "No" / "Do Not Execute
#      this Script." to "Yes" / "Execute this Script.", so the Manulife
#      script is now wired up to execute too (Springboot automation runner
#      picks up both rows now).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leading TestCaseID/TestCaseName columns; Execute_Flag,
# AutomationScriptName and Remarks shift left into A, B, C.
$ws.Range("A:B").Delete() | Out-Null

# Second data row (ManulifeLaunchURLTest) should now execute as well.
$ws.Range("A3").Value = "Yes"
$ws.Range("C3").Value = "Execute this Script."

# Leave the selection where the saved workbook had it.
$ws.Range("A3").Select() | Out-Null
